# Generate Report for Handoff
# A new handoff was produced for "b.md": update its status to "Ready for
# handoff" on the Overview sheet and on each locale sheet, and point the
# "Latest Handoff File" / "Latest Handoff Datetime" columns at the newly
# generated handoff file for that row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the "b.md" row. Columns B (zh-cn) and C (de-de)
# show the aggregate status for that file in each locale.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the "b.md" row.
#   B3 = Status
#   C3 = Latest Handoff File (hyperlinked)
#   D3 = Latest Handoff Datetime
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"

$zhcnC3 = $zhcn.Range("C3")
$zhcnC3.Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
if ($zhcnC3.Hyperlinks.Count -gt 0) {
    $zhcnC3.Hyperlinks.Item(1).TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
}

$zhcn.Range("D3").Value = "2016-03-09 12:56:37"

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the "b.md" row.
#   B3 = Status
#   C3 = Latest Handoff File (hyperlinked)
#   D3 = Latest Handoff Datetime
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"

$dedeC3 = $dede.Range("C3")
$dedeC3.Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
if ($dedeC3.Hyperlinks.Count -gt 0) {
    $dedeC3.Hyperlinks.Item(1).TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
}

$dede.Range("D3").Value = "2016-03-09 12:56:47"
